$d = $word.ActiveDocument

function Get-ParagraphStartingWith($doc, [string]$prefix) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    throw "No paragraph starting with '$prefix' was found"
}

# ---------------------------------------------------------------------------
# 1) "Currently working as Lead Developer & Scrum Master (Java, WSO2, Docker,
#    REST, OAuth, Ansible, Linux) React, Agile, Full-stack, Continuous
#    Integration & Delivery"
#    ->
#    "Currently working as Lead Developer & Scrum Master (Java, WSO2, Docker,
#    REST, OAuth, Ansible, Linux, JIRA, Git, Jenkins)"
#
#    Delete from "OAuth, Ansible, Linux)" through the end of the paragraph,
#    then re-insert the corrected tail as a new run appended at the end of
#    the paragraph. This leaves the earlier, untouched runs (including the
#    lone space-only run right before "OAuth") exactly as they were.
# ---------------------------------------------------------------------------
$p1 = Get-ParagraphStartingWith $d "Currently working"
$r1 = $p1.Range
$full1 = $r1.Text
$idxCut1 = $full1.IndexOf("OAuth, Ansible, Linux)")
$cutStart1 = $r1.Start + $idxCut1
$cutEnd1 = $r1.End - 1
$d.Range($cutStart1, $cutEnd1).Delete()

$p1b = Get-ParagraphStartingWith $d "Currently working"
$r1b = $p1b.Range
$tailPos1 = $r1b.End - 1
$d.Range($tailPos1, $tailPos1).InsertAfter("OAuth, Ansible, Linux, JIRA, Git, Jenkins)")

# ---------------------------------------------------------------------------
# 2) Remove the whole bullet paragraph "Continously improving the delivery
#    pipeline" (leftover from an incomplete copy & paste).
# ---------------------------------------------------------------------------
$pRemove = Get-ParagraphStartingWith $d "Continously improving"
$pRemove.Range.Delete()

# ---------------------------------------------------------------------------
# 3) "Open to new development challenges involving Java, Spring Boot, Docker,
#    API's, Microservices, CI/CD"
#    ->
#    "Open to new development challenges involving Java, Spring Boot, Docker,
#    API's, Microservices, React, Agile, Full-stack, Continuous Integration &
#    Delivery"
# ---------------------------------------------------------------------------
$p3 = Get-ParagraphStartingWith $d "Open to new development"
$r3 = $p3.Range
$full3 = $r3.Text
$idxCut3 = $full3.IndexOf("Microservices, CI/CD")
$cutStart3 = $r3.Start + $idxCut3
$cutEnd3 = $r3.End - 1
$d.Range($cutStart3, $cutEnd3).Delete()

$p3b = Get-ParagraphStartingWith $d "Open to new development"
$r3b = $p3b.Range
$tailPos3 = $r3b.End - 1
$d.Range($tailPos3, $tailPos3).InsertAfter("Microservices, React, Agile, Full-stack, Continuous Integration & Delivery")

Write-Host "Paragraph (Currently working): [$((Get-ParagraphStartingWith $d "Currently working").Range.Text)]"
Write-Host "Paragraph (Open to new development): [$((Get-ParagraphStartingWith $d "Open to new development").Range.Text)]"
